# "Generate Report for Handback"
#
# The handback-status report has one row per handed-off/handed-back file
# on each language sheet. For the 83fbcf0b... source file (row 3) on both
# the "zh-cn" and "de-de" sheets, record the handoff/handback report that
# was just regenerated by refreshing the "Correspond Handoff Datetime"
# (column E) and "Correspond Handback DateTime" (column H) timestamps.

$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets("zh-cn")
$wsDe = $wb.Worksheets("de-de")

# zh-cn sheet, row 3 -> 83fbcf0b-03ef-40b7-bbf4-d3452d6dd0cb... file
$wsZh.Range("E3").Value = "2016-03-19 03:38:20"
$wsZh.Range("H3").Value = "2016-03-19 03:39:04"

# de-de sheet, row 3 -> 83fbcf0b-03ef-40b7-bbf4-d3452d6dd0cb... file
$wsDe.Range("E3").Value = "2016-03-19 03:38:30"
$wsDe.Range("H3").Value = "2016-03-19 03:39:17"
